$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text shared strings with uniform run formatting) ---
$ws.Range("A8").Value = "Volume 32   Number  5"
$ws.Range("C9").Value = "Report Covering the Week  1/27/2025  Through  2/2/2025"

# --- Weekly crime statistics table updates (rows 14-33) ---
$ws.Range("G14").Value = "0"
$ws.Range("H14").Value = "***.*"
$ws.Range("I15").Value = 3
$ws.Range("N15").Value = 50
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 6
$ws.Range("H16").Value = -14.285714285714
$ws.Range("I16").Value = 8
$ws.Range("J16").Value = 9
$ws.Range("K16").Value = -11.111111111111
$ws.Range("L16").Value = -42.857142857142
$ws.Range("M16").Value = -50
$ws.Range("N16").Value = -86.206896551724
$ws.Range("C17").Value = 5
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 5
$ws.Range("H17").Value = 180
$ws.Range("I17").Value = 15
$ws.Range("K17").Value = 114.285714285714
$ws.Range("L17").Value = -21.052631578947
$ws.Range("M17").Value = 7.142857142857
$ws.Range("N17").Value = -6.25
$ws.Range("G18").Value = 2
$ws.Range("H18").Value = -50
$ws.Range("L18").Value = -66.666666666666
$ws.Range("M18").Value = -96.153846153846
$ws.Range("N18").Value = -98.461538461538
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = -80
$ws.Range("F19").Value = 10
$ws.Range("G19").Value = 21
$ws.Range("H19").Value = -52.380952380952
$ws.Range("J19").Value = 29
$ws.Range("K19").Value = -51.724137931034
$ws.Range("L19").Value = -44
$ws.Range("M19").Value = -22.222222222222
$ws.Range("N19").Value = -58.823529411764
$ws.Range("C20").Value = 7
$ws.Range("E20").Value = 75
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 16
$ws.Range("H20").Value = -43.75
$ws.Range("I20").Value = 11
$ws.Range("J20").Value = 18
$ws.Range("K20").Value = -38.888888888888
$ws.Range("L20").Value = 37.5
$ws.Range("M20").Value = -31.25
$ws.Range("N20").Value = -94.761904761904
$ws.Range("C21").Value = 14
$ws.Range("E21").Value = 27.272727272727
$ws.Range("F21").Value = 42
$ws.Range("G21").Value = 51
$ws.Range("H21").Value = -17.647058823529
$ws.Range("I21").Value = 52
$ws.Range("J21").Value = 67
$ws.Range("K21").Value = -22.388059701492
$ws.Range("L21").Value = -25.714285714285
$ws.Range("M21").Value = -42.222222222222
$ws.Range("N21").Value = -86.493506493506
$ws.Range("J22").Value = 3
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = "0"
$ws.Range("E23").Value = "***.*"
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = -11.111111111111
$ws.Range("I23").Value = 9
$ws.Range("K23").Value = -30.769230769230
$ws.Range("L23").Value = 28.571428571428
$ws.Range("M23").Value = 350
$ws.Range("C24").Value = 11
$ws.Range("D24").Value = 13
$ws.Range("E24").Value = -15.384615384615
$ws.Range("F24").Value = 32
$ws.Range("G24").Value = 49
$ws.Range("H24").Value = -34.693877551020
$ws.Range("I24").Value = 34
$ws.Range("J24").Value = 55
$ws.Range("K24").Value = -38.181818181818
$ws.Range("L24").Value = -45.161290322580
$ws.Range("M24").Value = -17.073170731707
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -25
$ws.Range("G25").Value = 12
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 13
$ws.Range("J25").Value = 14
$ws.Range("K25").Value = -7.142857142857
$ws.Range("L25").Value = 0
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = -20
$ws.Range("F26").Value = 23
$ws.Range("G26").Value = 28
$ws.Range("H26").Value = -17.857142857142
$ws.Range("I26").Value = 26
$ws.Range("J26").Value = 33
$ws.Range("K26").Value = -21.212121212121
$ws.Range("L26").Value = 4
$ws.Range("M26").Value = -29.729729729729
$ws.Range("F27").Value = 2
$ws.Range("I27").Value = 4
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = "0"
$ws.Range("E28").Value = "***.*"
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 50
$ws.Range("I28").Value = 4
$ws.Range("K28").Value = -20
$ws.Range("L28").Value = 300
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = -100
$ws.Range("F29").Value = 1
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 2
$ws.Range("K29").Value = -50
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = -100
$ws.Range("F30").Value = 1
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 2
$ws.Range("K30").Value = -50
$ws.Range("D33").Value = 2
$ws.Range("E33").Value = -100
$ws.Range("G33").Value = 2
$ws.Range("H33").Value = -50
$ws.Range("J33").Value = 2
$ws.Range("K33").Value = -50
